# Update the "Förändrad" date column (C2:C6) from 45224 (2023-10-25) to
# 45233 (2023-11-03) for each of the logging notification rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 45233
$ws.Range("C3").Value = 45233
$ws.Range("C4").Value = 45233
$ws.Range("C5").Value = 45233
$ws.Range("C6").Value = 45233
